$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Row 2 Col 4: "2321 (3.4 ± 0.59)" -> "2317 (3.4 ± 0.59)"
$table.Cell(2, 4).Range.Text = "2317 (3.4 ± 0.59)"

# Row 2 Col 5: "0.0033 (-2.5 ± 1.3)" -> "0.0032 (-2.5 ± 1.3)"
$table.Cell(2, 5).Range.Text = "0.0032 (-2.5 ± 1.3)"

# Row 2 Col 6: "208 (2.4 ± 1.3)" -> "214 (2.4 ± 1.3)"
$table.Cell(2, 6).Range.Text = "214 (2.4 ± 1.3)"

# Row 4 Col 4: "252 (2.4 ± 0.34)" -> "251 (2.4 ± 0.34)"
$table.Cell(4, 4).Range.Text = "251 (2.4 ± 0.34)"

# Row 4 Col 5: "0.011 (-2 ± 1)" -> "0.011 (-1.9 ± 1)"
$table.Cell(4, 5).Range.Text = "0.011 (-1.9 ± 1)"

# Row 4 Col 6: "64 (1.8 ± 1)" -> "63 (1.8 ± 1)"
$table.Cell(4, 6).Range.Text = "63 (1.8 ± 1)"

# Row 5 Col 4: "156 (2.2 ± 0.35)" -> "158 (2.2 ± 0.35)"
$table.Cell(5, 4).Range.Text = "158 (2.2 ± 0.35)"

# Row 5 Col 6: "48 (1.7 ± 1)" -> "49 (1.7 ± 1)"
$table.Cell(5, 6).Range.Text = "49 (1.7 ± 1)"

# Row 6 Col 3: "2.4 (0.38 ± 0.35)" -> "2.4 (0.37 ± 0.35)"
$table.Cell(6, 3).Range.Text = "2.4 (0.37 ± 0.35)"

# Row 7 Col 3: "2.8 (0.45 ± 0.36)" -> "2.7 (0.44 ± 0.36)"
$table.Cell(7, 3).Range.Text = "2.7 (0.44 ± 0.36)"

# Row 7 Col 4: "130 (2.1 ± 0.36)" -> "133 (2.1 ± 0.36)"
$table.Cell(7, 4).Range.Text = "133 (2.1 ± 0.36)"

# Row 7 Col 6: "45 (1.6 ± 1.1)" -> "46 (1.6 ± 1.1)"
$table.Cell(7, 6).Range.Text = "46 (1.6 ± 1.1)"

# Row 8 Col 4: "103 (2 ± 0.35)" -> "104 (2 ± 0.35)"
$table.Cell(8, 4).Range.Text = "104 (2 ± 0.35)"

# Row 9 Col 4: "168 (2.2 ± 0.33)" -> "169 (2.2 ± 0.33)"
$table.Cell(9, 4).Range.Text = "169 (2.2 ± 0.33)"

# Row 9 Col 5: "0.014 (-1.9 ± 1)" -> "0.013 (-1.9 ± 1)"
$table.Cell(9, 5).Range.Text = "0.013 (-1.9 ± 1)"

# Row 9 Col 6: "51 (1.7 ± 1)" -> "52 (1.7 ± 1)"
$table.Cell(9, 6).Range.Text = "52 (1.7 ± 1)"

# Row 10 Col 4: "118 (2.1 ± 0.36)" -> "119 (2.1 ± 0.36)"
$table.Cell(10, 4).Range.Text = "119 (2.1 ± 0.36)"

# Row 11 Col 3: "2.4 (0.39 ± 0.36)" -> "2.4 (0.38 ± 0.36)"
$table.Cell(11, 3).Range.Text = "2.4 (0.38 ± 0.36)"

# Row 12 Col 3: "3.5 (0.55 ± 0.35)" -> "3.5 (0.54 ± 0.35)"
$table.Cell(12, 3).Range.Text = "3.5 (0.54 ± 0.35)"

# Row 12 Col 4: "104 (2 ± 0.35)" -> "105 (2 ± 0.35)"
$table.Cell(12, 4).Range.Text = "105 (2 ± 0.35)"

# Row 12 Col 5: "0.018 (-1.7 ± 1.1)" -> "0.017 (-1.7 ± 1.1)"
$table.Cell(12, 5).Range.Text = "0.017 (-1.7 ± 1.1)"

# Row 12 Col 6: "39 (1.6 ± 1.1)" -> "40 (1.6 ± 1.1)"
$table.Cell(12, 6).Range.Text = "40 (1.6 ± 1.1)"

# Row 13 Col 5: "0.014 (-1.8 ± 1)" -> "0.014 (-1.9 ± 1)"
$table.Cell(13, 5).Range.Text = "0.014 (-1.9 ± 1)"

# Row 14 Col 6: "46 (1.6 ± 1.1)" -> "47 (1.7 ± 1.1)"
$table.Cell(14, 6).Range.Text = "47 (1.7 ± 1.1)"

# Row 15 Col 6: "52 (1.7 ± 1)" -> "51 (1.7 ± 1)"
$table.Cell(15, 6).Range.Text = "51 (1.7 ± 1)"

# Row 17 Col 4: "141 (2.1 ± 0.34)" -> "142 (2.1 ± 0.34)"
$table.Cell(17, 4).Range.Text = "142 (2.1 ± 0.34)"

# Row 17 Col 6: "47 (1.7 ± 1)" -> "46 (1.7 ± 1)"
$table.Cell(17, 6).Range.Text = "46 (1.7 ± 1)"

# Row 18 Col 3: "3 (0.48 ± 0.35)" -> "3 (0.47 ± 0.35)"
$table.Cell(18, 3).Range.Text = "3 (0.47 ± 0.35)"

# Row 18 Col 4: "122 (2.1 ± 0.35)" -> "123 (2.1 ± 0.35)"
$table.Cell(18, 4).Range.Text = "123 (2.1 ± 0.35)"

# Row 19 Col 6: "65 (1.8 ± 1)" -> "66 (1.8 ± 1)"
$table.Cell(19, 6).Range.Text = "66 (1.8 ± 1)"

# Row 20 Col 6: "67 (1.8 ± 1)" -> "68 (1.8 ± 1)"
$table.Cell(20, 6).Range.Text = "68 (1.8 ± 1)"

# Row 22 Col 4: "89 (1.9 ± 0.35)" -> "90 (2 ± 0.35)"
$table.Cell(22, 4).Range.Text = "90 (2 ± 0.35)"

# Row 23 Col 3: "1.8 (0.27 ± 0.34)" -> "1.8 (0.26 ± 0.34)"
$table.Cell(23, 3).Range.Text = "1.8 (0.26 ± 0.34)"

# Row 23 Col 4: "198 (2.3 ± 0.34)" -> "200 (2.3 ± 0.34)"
$table.Cell(23, 4).Range.Text = "200 (2.3 ± 0.34)"

# Row 23 Col 6: "56 (1.7 ± 1)" -> "57 (1.7 ± 1)"
$table.Cell(23, 6).Range.Text = "57 (1.7 ± 1)"

# Row 24 Col 3: "3.2 (0.51 ± 0.35)" -> "3.2 (0.51 ± 0.36)"
$table.Cell(24, 3).Range.Text = "3.2 (0.51 ± 0.36)"

# Row 24 Col 4: "113 (2 ± 0.35)" -> "113 (2.1 ± 0.36)"
$table.Cell(24, 4).Range.Text = "113 (2.1 ± 0.36)"

# Row 24 Col 6: "41 (1.6 ± 1.1)" -> "42 (1.6 ± 1.1)"
$table.Cell(24, 6).Range.Text = "42 (1.6 ± 1.1)"

# Row 25 Col 3: "3.5 (0.54 ± 0.37)" -> "3.4 (0.54 ± 0.36)"
$table.Cell(25, 3).Range.Text = "3.4 (0.54 ± 0.36)"

# Row 25 Col 4: "106 (2 ± 0.37)" -> "107 (2 ± 0.36)"
$table.Cell(25, 4).Range.Text = "107 (2 ± 0.36)"

# Row 25 Col 5: "0.018 (-1.7 ± 1.1)" -> "0.017 (-1.7 ± 1.1)"
$table.Cell(25, 5).Range.Text = "0.017 (-1.7 ± 1.1)"

# Row 26 Col 3: "5.5 (0.74 ± 0.34)" -> "5.5 (0.74 ± 0.33)"
$table.Cell(26, 3).Range.Text = "5.5 (0.74 ± 0.33)"

# Row 26 Col 4: "66 (1.8 ± 0.34)" -> "67 (1.8 ± 0.33)"
$table.Cell(26, 4).Range.Text = "67 (1.8 ± 0.33)"

# Row 26 Col 6: "30 (1.5 ± 1.1)" -> "31 (1.5 ± 1.1)"
$table.Cell(26, 6).Range.Text = "31 (1.5 ± 1.1)"

# Row 27 Col 3: "3.8 (0.58 ± 0.35)" -> "3.8 (0.57 ± 0.35)"
$table.Cell(27, 3).Range.Text = "3.8 (0.57 ± 0.35)"

# Row 27 Col 5: "0.019 (-1.7 ± 1.1)" -> "0.018 (-1.7 ± 1.1)"
$table.Cell(27, 5).Range.Text = "0.018 (-1.7 ± 1.1)"

# Row 27 Col 6: "37 (1.6 ± 1.1)" -> "38 (1.6 ± 1.1)"
$table.Cell(27, 6).Range.Text = "38 (1.6 ± 1.1)"

# Row 29 Col 3: "3.6 (0.56 ± 0.37)" -> "3.6 (0.55 ± 0.36)"
$table.Cell(29, 3).Range.Text = "3.6 (0.55 ± 0.36)"

# Row 29 Col 4: "101 (2 ± 0.37)" -> "102 (2 ± 0.36)"
$table.Cell(29, 4).Range.Text = "102 (2 ± 0.36)"

# Row 30 Col 3: "2.4 (0.39 ± 0.35)" -> "2.4 (0.38 ± 0.35)"
$table.Cell(30, 3).Range.Text = "2.4 (0.38 ± 0.35)"

# Row 30 Col 4: "151 (2.2 ± 0.35)" -> "152 (2.2 ± 0.35)"
$table.Cell(30, 4).Range.Text = "152 (2.2 ± 0.35)"

# Row 30 Col 5: "0.014 (-1.8 ± 1)" -> "0.015 (-1.8 ± 1)"
$table.Cell(30, 5).Range.Text = "0.015 (-1.8 ± 1)"

# Row 30 Col 6: "49 (1.7 ± 1)" -> "48 (1.7 ± 1)"
$table.Cell(30, 6).Range.Text = "48 (1.7 ± 1)"

# Row 31 Col 3: "3.4 (0.54 ± 0.35)" -> "3.5 (0.54 ± 0.35)"
$table.Cell(31, 3).Range.Text = "3.5 (0.54 ± 0.35)"

# Row 32 Col 4: "116 (2.1 ± 0.36)" -> "117 (2.1 ± 0.36)"
$table.Cell(32, 4).Range.Text = "117 (2.1 ± 0.36)"

# Row 32 Col 6: "42 (1.6 ± 1.1)" -> "41 (1.6 ± 1.1)"
$table.Cell(32, 6).Range.Text = "41 (1.6 ± 1.1)"

# Row 33 Col 4: "114 (2.1 ± 0.36)" -> "115 (2.1 ± 0.36)"
$table.Cell(33, 4).Range.Text = "115 (2.1 ± 0.36)"

# Row 33 Col 5: "0.017 (-1.8 ± 1.1)" -> "0.016 (-1.8 ± 1.1)"
$table.Cell(33, 5).Range.Text = "0.016 (-1.8 ± 1.1)"

# Row 34 Col 6: "47 (1.7 ± 1)" -> "48 (1.7 ± 1)"
$table.Cell(34, 6).Range.Text = "48 (1.7 ± 1)"

# Row 35 Col 3: "3.4 (0.53 ± 0.34)" -> "3.3 (0.53 ± 0.34)"
$table.Cell(35, 3).Range.Text = "3.3 (0.53 ± 0.34)"

# Row 35 Col 4: "108 (2 ± 0.34)" -> "109 (2 ± 0.34)"
$table.Cell(35, 4).Range.Text = "109 (2 ± 0.34)"

# Row 35 Col 5: "0.017 (-1.7 ± 1.1)" -> "0.017 (-1.8 ± 1.1)"
$table.Cell(35, 5).Range.Text = "0.017 (-1.8 ± 1.1)"

# Row 35 Col 6: "40 (1.6 ± 1.1)" -> "41 (1.6 ± 1.1)"
$table.Cell(35, 6).Range.Text = "41 (1.6 ± 1.1)"

# Row 36 Col 5: "0.017 (-1.7 ± 1.1)" -> "0.017 (-1.8 ± 1.1)"
$table.Cell(36, 5).Range.Text = "0.017 (-1.8 ± 1.1)"
